# Fruta / hortaliza, semanal
# Insert a new weekly record as row 131 in the Jengibre sheet, pushing all
# subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 131 (shifts rows 131-145 down to 132-146)
$ws.Rows("131:131").Insert()

# Populate the new row 131 with the new weekly data point
$ws.Cells.Item(131, 1).Value = 9
$ws.Cells.Item(131, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(131, 3).Value = "Metropolitana"
$ws.Cells.Item(131, 4).Value = 45194
$ws.Cells.Item(131, 5).Value = 13
$ws.Cells.Item(131, 6).Value = 100114007
$ws.Cells.Item(131, 7).Value = "Jengibre"
$ws.Cells.Item(131, 8).Value = "Sin especificar"
$ws.Cells.Item(131, 9).Value = "Primera"
$ws.Cells.Item(131, 10).Value = 340
$ws.Cells.Item(131, 11).Value = 20000
$ws.Cells.Item(131, 12).Value = 21000
$ws.Cells.Item(131, 13).Value = 20500
$ws.Cells.Item(131, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(131, 15).Value = "Perú"
$ws.Cells.Item(131, 16).Value = 1577
$ws.Cells.Item(131, 17).Value = 13
$ws.Cells.Item(131, 18).Value = "Hortaliza"
